$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue "D2" "59.912.43"
Set-TextValue "E2" "  +1.35%  "
Set-TextValue "D3" "2.307.16"
Set-TextValue "E3" "  -0.15%  "
Set-TextValue "E4" "  -0.04%  "
Set-TextValue "D5" "541.55"
Set-TextValue "E5" "  -0.12%  "
Set-TextValue "D6" "129.57"
Set-TextValue "E6" "  -1.99%  "
Set-TextValue "E7" "  -0.03%  "
Set-TextValue "D8" "0.573"
Set-TextValue "E8" "  -2.37%  "
Set-TextValue "D9" "2.305.52"
Set-TextValue "E9" "  -0.16%  "
Set-TextValue "E10" "  -0.17%  "
Set-TextValue "E11" "  +1.95%  "
Set-TextValue "E12" "  -0.26%  "
Set-TextValue "E13" "  +0.16%  "
Set-TextValue "D14" "23.33"
Set-TextValue "E14" "  -2.15%  "
Set-TextValue "D15" "59.907.75"
Set-TextValue "E15" "  +1.54%  "
Set-TextValue "D16" "2.718.55"
Set-TextValue "E16" "  -0.19%  "
Set-TextValue "E17" "  -0.93%  "
Set-TextValue "D18" "2.301.68"
Set-TextValue "E18" "  -0.49%  "
Set-TextValue "D19" "10.49"
Set-TextValue "E19" "  -1.08%  "
Set-TextValue "E20" "  -2.07%  "
Set-TextValue "D21" "312.56"
Set-TextValue "E21" "  +0.02%  "
Set-TextValue "E22" "  -0.29%  "
Set-TextValue "E23" "  -0.15%  "
Set-TextValue "D24" "63.65"
Set-TextValue "D25" "0.170"
Set-TextValue "E25" "  -2.02%  "
Set-TextValue "E26" "  -0.02%  "
Set-TextValue "E27" "  -2.78%  "
Set-TextValue "E28" "  +4.15%  "
Set-TextValue "D29" "171.49"
Set-TextValue "E29" "  +0.82%  "
Set-TextValue "E30" "  -0.66%  "
Set-TextValue "E31" "  +0.07%  "
Set-TextValue "D32" "0.0₃0726"
Set-TextValue "E32" "  -1.79%  "
Set-TextValue "E33" "  -0.75%  "
Set-TextValue "E34" "  +3.33%  "
Set-TextValue "E35" "  -1.31%  "
Set-TextValue "D37" "17.69"
Set-TextValue "E37" "  -0.83%  "
Set-TextValue "E38" "  -0.03%  "
Set-TextValue "D39" "4.00"
Set-TextValue "E39" "  -1.60%  "
Set-TextValue "D40" "316.21"
Set-TextValue "E40" "  +2.53%  "
Set-TextValue "E41" "  +0.65%  "
Set-TextValue "E42" "  -0.51%  "
Set-TextValue "D43" "136.06"
Set-TextValue "E43" "  -3.58%  "
Set-TextValue "E44" "  -0.42%  "
Set-TextValue "E45" "  -1.95%  "
Set-TextValue "E46" "  +2.17%  "
Set-TextValue "E47" "  +2.57%  "
Set-TextValue "D48" "0.0491"
Set-TextValue "E48" "  -0.91%  "
Set-TextValue "D49" "0.0₆0225"
Set-TextValue "E49" "  +22.25%  "
Set-TextValue "E50" "  +1.26%  "
Set-TextValue "E51" "  +0.06%  "
